# Regenerate save_data to use K (strikeouts) instead of Strike# (total
# strikes thrown) in column G, and write the recalculated s_vals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G ("K") values, keyed by row number (row 2 is the first data
# row, row 34 is the last).
$kValues = @{
    2  = 7
    3  = 6
    4  = 2
    5  = 4
    6  = 4
    7  = 0
    8  = 4
    9  = 8
    10 = 2
    11 = 5
    12 = 4
    13 = 3
    14 = 4
    15 = 4
    16 = 8
    17 = 4
    18 = 7
    19 = 4
    20 = 5
    21 = 2
    22 = 2
    23 = 5
    24 = 10
    25 = 1
    26 = 3
    27 = 3
    28 = 3
    29 = 3
    30 = 3
    31 = 4
    32 = 4
    33 = 2
    34 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
